$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("C4").Value = "广州·异世界夜宴（取消）"
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F6").Value = 364
$ws1.Range("G6").Value = 41
$ws1.Range("F8").Value = 148
$ws1.Range("F9").Value = 241
$ws1.Range("F10").Value = 218
$ws1.Range("F11").Value = 5939
$ws1.Range("F12").Value = 55
$ws1.Range("F13").Value = 46
$ws1.Range("F14").Value = 492
$ws1.Range("F17").Value = 356
$ws1.Range("F21").Value = 706
$ws1.Range("F22").Value = 138
$ws1.Range("F23").Value = 94
$ws1.Range("F24").Value = 310
$ws1.Range("F27").Value = 1808
$ws1.Range("F28").Value = 467

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 266
$ws2.Range("F5").Value = 269
$ws2.Range("F6").Value = 301

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 575
$ws4.Range("C5").Value = "广州·异世界夜宴（取消）"
$ws4.Range("G5").Value = "不可售"
$ws4.Range("F8").Value = 364
$ws4.Range("G8").Value = 41
$ws4.Range("F10").Value = 148
$ws4.Range("F11").Value = 241
$ws4.Range("F12").Value = 218
$ws4.Range("F13").Value = 5939
$ws4.Range("F14").Value = 55
$ws4.Range("F15").Value = 46
$ws4.Range("F16").Value = 266
$ws4.Range("F17").Value = 492
$ws4.Range("F20").Value = 356
$ws4.Range("F25").Value = 269
$ws4.Range("F26").Value = 301
$ws4.Range("F28").Value = 706
$ws4.Range("F32").Value = 138
$ws4.Range("F33").Value = 94
$ws4.Range("F34").Value = 310
$ws4.Range("F37").Value = 1808
$ws4.Range("F38").Value = 467
